# Sigi-9: Generate filename and sequence number from seal id
#
# A new seal (column D, Seal ID 115) is added to the comparison sheet, and
# the manually-maintained "Filename" / "Sequence" rows are removed entirely
# -- those values are now generated programmatically from the Seal ID
# instead of being hand-kept in the spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new seal's data in column D ------------------------------
$ws.Range("D1").Value = 115                     # SEAL ID
$ws.Range("D2").Value = "Seal"                   # TYPE
$ws.Range("D57").Value = "Another great title"   # Title
$ws.Range("D58").Value = "Clark"                 # Editor forename
$ws.Range("D59").Value = "Kent"                  # Editor surname
$ws.Range("D45").Value = "Another Edition"       # EDITION(S)

# --- Remove the now-obsolete "Filename" (row 60) and "Sequence" (row 61)
$ws.Range("A60:A61").EntireRow.Delete()

# --- Leave the selection where the edit happened (cosmetic, matches the
#     author's saved cursor position) -----------------------------------
[void]$ws.Range("A60:XFD60").Select()
